# Insert a new data row at row 59 (pushing the existing rows 59-82 down to 60-83)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(59).Insert()

$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value = 44609
$ws.Cells.Item(59, 5).Value = 15
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100102
$ws.Cells.Item(59, 8).Value = "Cítricos"
$ws.Cells.Item(59, 9).Value = 100102005
$ws.Cells.Item(59, 10).Value = "Naranja"
$ws.Cells.Item(59, 11).Value = "Valencia"
$ws.Cells.Item(59, 12).Value = "Segunda"
$ws.Cells.Item(59, 13).Value = 250
$ws.Cells.Item(59, 14).Value = 850
$ws.Cells.Item(59, 15).Value = 900
$ws.Cells.Item(59, 16).Value = 875
$ws.Cells.Item(59, 17).Value = '$/kilo (en caja de 20 kilos)'
$ws.Cells.Item(59, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(59, 19).Value = 875
$ws.Cells.Item(59, 20).Value = 1
